$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04332866666666666
$ws.Range("H2").Value = 0.129986
$ws.Range("I2").Value = 0.2246397599897691
$ws.Range("J2").Value = 0.2246397599897691
$ws.Range("M2").Value = 0.2799683333333333
$ws.Range("N2").Value = 0.839905
$ws.Range("O2").Value = 0.0294305463214559
$ws.Range("P2").Value = 0.0294305463214559
$ws.Range("Q2").Value = 0.01213065459222222
$ws.Range("R2").Value = 0.10917589133
$ws.Range("S2").Value = 0.006611270862019635
$ws.Range("T2").Value = 0.006611270862019637
$ws.Range("G3").Value = 0.04332866666666666
$ws.Range("H3").Value = 0.129986
$ws.Range("I3").Value = 0.2246397599897691
$ws.Range("J3").Value = 0.2246397599897691
$ws.Range("O3").Value = 0.2486942046732164
$ws.Range("P3").Value = 0.2486942046732163
$ws.Range("Q3").Value = 0.1025065407562222
$ws.Range("R3").Value = 0.922558866806
$ws.Range("S3").Value = 0.05586660644863784
$ws.Range("T3").Value = 0.05586660644863785
$ws.Range("G4").Value = 0.04332866666666666
$ws.Range("H4").Value = 0.129986
$ws.Range("I4").Value = 0.2246397599897691
$ws.Range("J4").Value = 0.2246397599897691
$ws.Range("M4").Value = 6.86709
$ws.Range("N4").Value = 20.60127
$ws.Range("O4").Value = 0.7218752490053277
$ws.Range("P4").Value = 0.7218752490053277
$ws.Range("Q4").Value = 0.29754185358
$ws.Range("R4").Value = 2.67787668222
$ws.Range("S4").Value = 0.1621618826791116
$ws.Range("T4").Value = 0.1621618826791117
$ws.Range("I5").Value = 0.5955530362469368
$ws.Range("J5").Value = 0.5955530362469369
$ws.Range("M5").Value = 0.2799683333333333
$ws.Range("N5").Value = 0.839905
$ws.Range("O5").Value = 0.0294305463214559
$ws.Range("P5").Value = 0.0294305463214559
$ws.Range("Q5").Value = 0.03216014909555556
$ws.Range("R5").Value = 0.28944134186
$ws.Range("S5").Value = 0.01752745122014918
$ws.Range("T5").Value = 0.01752745122014918
$ws.Range("I6").Value = 0.5955530362469368
$ws.Range("J6").Value = 0.5955530362469369
$ws.Range("O6").Value = 0.2486942046732164
$ws.Range("P6").Value = 0.2486942046732163
$ws.Range("S6").Value = 0.1481105886901511
$ws.Range("T6").Value = 0.1481105886901511
$ws.Range("I7").Value = 0.5955530362469368
$ws.Range("J7").Value = 0.5955530362469369
$ws.Range("M7").Value = 6.86709
$ws.Range("N7").Value = 20.60127
$ws.Range("O7").Value = 0.7218752490053277
$ws.Range("P7").Value = 0.7218752490053277
$ws.Range("Q7").Value = 0.78882720636
$ws.Range("R7").Value = 7.099444857240001
$ws.Range("S7").Value = 0.4299149963366365
$ws.Range("T7").Value = 0.4299149963366365
$ws.Range("G8").Value = 0.03468133333333333
$ws.Range("H8").Value = 0.104044
$ws.Range("I8").Value = 0.1798072037632941
$ws.Range("J8").Value = 0.1798072037632941
$ws.Range("M8").Value = 0.2799683333333333
$ws.Range("N8").Value = 0.839905
$ws.Range("O8").Value = 0.0294305463214559
$ws.Range("P8").Value = 0.0294305463214559
$ws.Range("Q8").Value = 0.009709675091111111
$ws.Range("R8").Value = 0.08738707581999999
$ws.Range("S8").Value = 0.005291824239287085
$ws.Range("T8").Value = 0.005291824239287086
$ws.Range("G9").Value = 0.03468133333333333
$ws.Range("H9").Value = 0.104044
$ws.Range("I9").Value = 0.1798072037632941
$ws.Range("J9").Value = 0.1798072037632941
$ws.Range("O9").Value = 0.2486942046732164
$ws.Range("P9").Value = 0.2486942046732163
$ws.Range("Q9").Value = 0.08204876314711111
$ws.Range("R9").Value = 0.738438868324
$ws.Range("S9").Value = 0.04471700953442737
$ws.Range("T9").Value = 0.04471700953442737
$ws.Range("G10").Value = 0.03468133333333333
$ws.Range("H10").Value = 0.104044
$ws.Range("I10").Value = 0.1798072037632941
$ws.Range("J10").Value = 0.1798072037632941
$ws.Range("M10").Value = 6.86709
$ws.Range("N10").Value = 20.60127
$ws.Range("O10").Value = 0.7218752490053277
$ws.Range("P10").Value = 0.7218752490053277
$ws.Range("Q10").Value = 0.23815983732
$ws.Range("R10").Value = 2.14343853588
$ws.Range("S10").Value = 0.1297983699895796
$ws.Range("T10").Value = 0.1297983699895796
